$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (they are stored as text in the source).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "44.760.26"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.260.33"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").Value = "301.53"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").Value = "94.20"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").Value = "34.02"
$ws.Range("E10").Value = "  -3.84%  "
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "2.602.93"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "2.255.41"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("D16").Value = "13.55"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  -5.62%  "
$ws.Range("D18").Value = "44.677.27"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "12.87"
$ws.Range("E19").Value = "  +6.60%  "
$ws.Range("D20").Value = "0.0₃0921"
$ws.Range("E20").Value = "  -3.30%  "
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").Value = "65.55"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "237.84"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").Value = "2.88"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  -5.42%  "
$ws.Range("D27").Value = "41.28"
$ws.Range("E27").Value = "  +9.53%  "
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("E29").Value = "  -3.25%  "
$ws.Range("D30").Value = "19.52"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "152.04"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").Value = "5.53"
$ws.Range("E32").Value = "  -9.38%  "
$ws.Range("D33").Value = "0.0789"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").Value = "2.95"
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("D38").Value = "1.74"
$ws.Range("E38").Value = "  -5.76%  "
$ws.Range("D39").Value = "4.01"
$ws.Range("E39").Value = "  +5.57%  "
$ws.Range("D40").Value = "0.0309"
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("D41").Value = "3.23"
$ws.Range("E41").Value = "  -5.50%  "
$ws.Range("D42").Value = "13.57"
$ws.Range("E42").Value = "  -10.14%  "
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "1.89"
$ws.Range("E44").Value = "  +7.43%  "
$ws.Range("D45").Value = "1.736.93"
$ws.Range("E45").Value = "  -6.25%  "
$ws.Range("D46").Value = "0.194"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "69.32"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "75.75"
$ws.Range("E48").Value = "  -5.42%  "
$ws.Range("D49").Value = "95.90"
$ws.Range("E49").Value = "  -3.59%  "
$ws.Range("D50").Value = "53.56"
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("D51").Value = "4.68"
$ws.Range("E51").Value = "  -4.74%  "

Write-Host "Updated cryptos list"
